# Applies the cryptos.xlsx price/volume refresh described in the commit
# 'Updated cryptos list on Thu May  2 09:45:32 UTC 2024 with GitHub Actions'.
# All Coin/Link/Price/Volume cells in this sheet are stored as plain text,
# so numeric-looking Price values are written via a Text-formatted cell
# (then restored to the Normal style) to stop Excel from auto-converting
# them into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.799.06"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.946.01"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.13%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +5.23%  "
$ws.Range("D9").Value = "2.937.46"
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("E10").Value = "  +3.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  +4.17%  "
$ws.Range("E13").Value = "  +5.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("D16").Value = "3.432.94"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("E17").Value = "  +8.87%  "
$ws.Range("D18").Value = "2.944.31"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").Value = "57.774.04"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "415.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.07%  "
$ws.Range("E22").Value = "  +7.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  +6.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0962"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.951"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.13%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.36%  "
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("E37").Value = "  +7.67%  "
$ws.Range("D38").Value = "0.0₃0697"
$ws.Range("E38").Value = "  +14.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +16.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "382.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.57%  "
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("D44").Value = "2.705.05"
$ws.Range("E44").Value = "  +4.18%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.86%  "
$ws.Range("E47").Value = "  +4.30%  "
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("E51").Value = "  +3.74%  "
